$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 966.41174
$ws.Cells.Item(28, 9).Value = 794.1667
$ws.Cells.Item(28, 10).Value = 1379.8
$ws.Cells.Item(28, 11).Value = 794.1667
$ws.Cells.Item(28, 12).Value = 1379.8
$ws.Cells.Item(28, 13).Value = -309.1667
$ws.Cells.Item(28, 14).Value = -2349.8
$ws.Cells.Item(32, 8).Value = 16671167
$ws.Cells.Item(32, 9).Value = 2000
$ws.Cells.Item(32, 11).Value = 2000
$ws.Cells.Item(32, 13).Value = -1674
$ws.Cells.Item(40, 8).Value = 5453.222
$ws.Cells.Item(40, 9).Value = 6318
$ws.Cells.Item(40, 10).Value = 3723.6667
$ws.Cells.Item(40, 11).Value = 6318
$ws.Cells.Item(40, 12).Value = 3723.6667
$ws.Cells.Item(40, 13).Value = -6143
$ws.Cells.Item(40, 14).Value = -4073.6667
$ws.Cells.Item(92, 8).Value = 454.77777
$ws.Cells.Item(92, 9).Value = 467.3913
$ws.Cells.Item(92, 11).Value = 467.3913
$ws.Cells.Item(92, 13).Value = 780.6087
$ws.Cells.Item(112, 8).Value = 1655.3
$ws.Cells.Item(112, 10).Value = 1802.6
$ws.Cells.Item(112, 12).Value = 5407.799999999999
$ws.Cells.Item(112, 14).Value = -7623.799999999999
$ws.Cells.Item(137, 8).Value = 15407848
$ws.Cells.Item(137, 9).Value = 25033978
$ws.Cells.Item(137, 10).Value = 6039.8
$ws.Cells.Item(137, 11).Value = 75101934
$ws.Cells.Item(137, 12).Value = 18119.4
$ws.Cells.Item(137, 13).Value = -75099384
$ws.Cells.Item(137, 14).Value = -23219.4
$ws.Cells.Item(138, 8).Value = 3441.75
$ws.Cells.Item(138, 9).Value = 2034.0625
$ws.Cells.Item(138, 10).Value = 4849.4375
$ws.Cells.Item(138, 11).Value = 6102.1875
$ws.Cells.Item(138, 12).Value = 14548.3125
$ws.Cells.Item(138, 13).Value = -962.1875
$ws.Cells.Item(138, 14).Value = -24828.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4819
$ws.Cells.Item(32, 9).Value = 5080.5947
$ws.Cells.Item(32, 11).Value = 5080.5947
$ws.Cells.Item(32, 13).Value = -4793.5947
$ws.Cells.Item(74, 8).Value = 4265.8667
$ws.Cells.Item(74, 9).Value = 3898.5
$ws.Cells.Item(74, 10).Value = 5000.6
$ws.Cells.Item(74, 11).Value = 3898.5
$ws.Cells.Item(74, 12).Value = 5000.6
$ws.Cells.Item(74, 13).Value = -3024.5
$ws.Cells.Item(74, 14).Value = -6748.6
$ws.Cells.Item(77, 8).Value = 4265.8667
$ws.Cells.Item(77, 9).Value = 3898.5
$ws.Cells.Item(77, 10).Value = 5000.6
$ws.Cells.Item(77, 11).Value = 19492.5
$ws.Cells.Item(77, 12).Value = 25003
$ws.Cells.Item(77, 13).Value = -15124.5
$ws.Cells.Item(77, 14).Value = -33739
$ws.Cells.Item(97, 8).Value = 4326.7896
$ws.Cells.Item(97, 9).Value = 3456.0557
$ws.Cells.Item(97, 11).Value = 3456.0557
$ws.Cells.Item(97, 13).Value = -2960.0557
$ws.Cells.Item(110, 8).Value = 3323.2246
$ws.Cells.Item(110, 9).Value = 3138.2285
$ws.Cells.Item(110, 11).Value = 3138.2285
$ws.Cells.Item(110, 13).Value = -1093.2285

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 4384.3
$ws.Cells.Item(94, 9).Value = 4004.125
$ws.Cells.Item(94, 11).Value = 4004.125
$ws.Cells.Item(94, 13).Value = -3553.125
$ws.Cells.Item(105, 8).Value = 5547.727
$ws.Cells.Item(105, 9).Value = 4732.143
$ws.Cells.Item(105, 11).Value = 4732.143
$ws.Cells.Item(105, 13).Value = -2985.143
$ws.Cells.Item(107, 8).Value = 3607.0833
$ws.Cells.Item(107, 9).Value = 2981.8948
$ws.Cells.Item(107, 10).Value = 5982.8
$ws.Cells.Item(107, 11).Value = 2981.8948
$ws.Cells.Item(107, 12).Value = 5982.8
$ws.Cells.Item(107, 13).Value = -1061.8948
$ws.Cells.Item(107, 14).Value = -9822.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6300.0454
$ws.Cells.Item(31, 9).Value = 9062.75
$ws.Cells.Item(31, 10).Value = 2984.8
$ws.Cells.Item(31, 11).Value = 9062.75
$ws.Cells.Item(31, 12).Value = 2984.8
$ws.Cells.Item(31, 13).Value = -8767.75
$ws.Cells.Item(31, 14).Value = -3574.8
$ws.Cells.Item(34, 8).Value = 6300.0454
$ws.Cells.Item(34, 9).Value = 9062.75
$ws.Cells.Item(34, 10).Value = 2984.8
$ws.Cells.Item(34, 11).Value = 9062.75
$ws.Cells.Item(34, 12).Value = 2984.8
$ws.Cells.Item(34, 13).Value = -8860.75
$ws.Cells.Item(34, 14).Value = -3388.8
$ws.Cells.Item(58, 8).Value = 2638.0588
$ws.Cells.Item(58, 9).Value = 1881.375
$ws.Cells.Item(58, 11).Value = 1881.375
$ws.Cells.Item(58, 13).Value = -1678.375
$ws.Cells.Item(105, 8).Value = 1603.6111
$ws.Cells.Item(105, 9).Value = 1539.125
$ws.Cells.Item(105, 11).Value = 1539.125
$ws.Cells.Item(105, 13).Value = 207.875
$ws.Cells.Item(107, 8).Value = 382.73914
$ws.Cells.Item(107, 10).Value = 704.8333
$ws.Cells.Item(107, 12).Value = 704.8333
$ws.Cells.Item(107, 14).Value = -4544.8333
$ws.Cells.Item(119, 8).Value = 37900
$ws.Cells.Item(119, 10).Value = 37900
$ws.Cells.Item(119, 12).Value = 37900
$ws.Cells.Item(119, 14).Value = -47576
$ws.Cells.Item(136, 8).Value = 2638.0588
$ws.Cells.Item(136, 9).Value = 1881.375
$ws.Cells.Item(136, 11).Value = 5644.125
$ws.Cells.Item(136, 13).Value = -3094.125
$ws.Cells.Item(141, 8).Value = 60000
$ws.Cells.Item(141, 10).Value = 60000
$ws.Cells.Item(141, 12).Value = 60000
$ws.Cells.Item(141, 14).Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 739.46155
$ws.Cells.Item(5, 9).Value = 1570.8334
$ws.Cells.Item(5, 10).Value = 490.05
$ws.Cells.Item(5, 11).Value = 4712.5002
$ws.Cells.Item(5, 12).Value = 1470.15
$ws.Cells.Item(5, 13).Value = -4600.5002
$ws.Cells.Item(5, 14).Value = -1694.15
$ws.Cells.Item(68, 8).Value = 1409.8334
$ws.Cells.Item(68, 10).Value = 1482.6666
$ws.Cells.Item(68, 12).Value = 4447.9998
$ws.Cells.Item(68, 14).Value = -6069.9998
$ws.Cells.Item(71, 8).Value = 1409.8334
$ws.Cells.Item(71, 10).Value = 1482.6666
$ws.Cells.Item(71, 12).Value = 13343.9994
$ws.Cells.Item(71, 14).Value = -21455.9994
$ws.Cells.Item(107, 8).Value = 1417.8667
$ws.Cells.Item(107, 9).Value = 877.1
$ws.Cells.Item(107, 10).Value = 2499.4
$ws.Cells.Item(107, 11).Value = 2631.3
$ws.Cells.Item(107, 12).Value = 7498.200000000001
$ws.Cells.Item(107, 13).Value = -711.3000000000002
$ws.Cells.Item(107, 14).Value = -11338.2
$ws.Cells.Item(113, 8).Value = 651
$ws.Cells.Item(113, 9).Value = 1085.6
$ws.Cells.Item(113, 10).Value = 523.17645
$ws.Cells.Item(113, 11).Value = 3256.8
$ws.Cells.Item(113, 12).Value = 1569.52935
$ws.Cells.Item(113, 13).Value = -1086.8
$ws.Cells.Item(113, 14).Value = -5909.529350000001
$ws.Cells.Item(116, 8).Value = 2859.8
$ws.Cells.Item(116, 9).Value = 400
$ws.Cells.Item(116, 11).Value = 1200
$ws.Cells.Item(116, 13).Value = 2242
$ws.Cells.Item(135, 8).Value = 739.46155
$ws.Cells.Item(135, 9).Value = 1570.8334
$ws.Cells.Item(135, 10).Value = 490.05
$ws.Cells.Item(135, 11).Value = 14137.5006
$ws.Cells.Item(135, 12).Value = 4410.45
$ws.Cells.Item(135, 13).Value = -11602.5006
$ws.Cells.Item(135, 14).Value = -9480.450000000001
$ws.Cells.Item(136, 8).Value = 3693.3333
$ws.Cells.Item(136, 9).Value = 2600
$ws.Cells.Item(136, 10).Value = 4090.9092
$ws.Cells.Item(136, 11).Value = 7800
$ws.Cells.Item(136, 12).Value = 12272.7276
$ws.Cells.Item(136, 13).Value = -2700
$ws.Cells.Item(136, 14).Value = -22472.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3126.1924
$ws.Cells.Item(102, 9).Value = 3074.2083
$ws.Cells.Item(102, 11).Value = 3074.2083
$ws.Cells.Item(102, 13).Value = -1452.2083
$ws.Cells.Item(113, 8).Value = 12441.571
$ws.Cells.Item(113, 9).Value = 9452.154
$ws.Cells.Item(113, 10).Value = 15032.4
$ws.Cells.Item(113, 11).Value = 9452.154
$ws.Cells.Item(113, 12).Value = 15032.4
$ws.Cells.Item(113, 13).Value = -7282.154
$ws.Cells.Item(113, 14).Value = -19372.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6278.3335
$ws.Cells.Item(40, 9).Value = 7431.5
$ws.Cells.Item(40, 11).Value = 7431.5
$ws.Cells.Item(40, 13).Value = -7295.5
$ws.Cells.Item(46, 8).Value = 4565.839
$ws.Cells.Item(46, 9).Value = 2837.5
$ws.Cells.Item(46, 11).Value = 2837.5
$ws.Cells.Item(46, 13).Value = -2649.5
$ws.Cells.Item(122, 8).Value = 6245.7144
$ws.Cells.Item(122, 10).Value = 7996.3335
$ws.Cells.Item(122, 12).Value = 23989.0005
$ws.Cells.Item(122, 14).Value = -28889.0005
$ws.Cells.Item(136, 8).Value = 4529.524
$ws.Cells.Item(136, 9).Value = 5840
$ws.Cells.Item(136, 11).Value = 17520
$ws.Cells.Item(136, 13).Value = -14970

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 475
$ws.Cells.Item(113, 9).Value = 552.2
$ws.Cells.Item(113, 10).Value = 217.66667
$ws.Cells.Item(113, 11).Value = 1656.6
$ws.Cells.Item(113, 12).Value = 653.00001
$ws.Cells.Item(113, 13).Value = 513.3999999999999
$ws.Cells.Item(113, 14).Value = -4993.00001
$ws.Cells.Item(136, 8).Value = 4076.5833
$ws.Cells.Item(136, 9).Value = 2609.72
$ws.Cells.Item(136, 11).Value = 7829.16
$ws.Cells.Item(136, 13).Value = -5279.16
